$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a genuinely numeric value to a cell that lives in a
# column formatted as Text (numFmtId 49 / "@"), so it is NOT coerced to a
# text string the way a bare `.Value = <number>` assignment would be.
function Set-NumericValue($cell, $value) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "General"
    $cell.Value = $value
    $cell.NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# 1. Swap the D/E columns for the existing "EindArtikel" rows (2-11):
#    D used to hold the numeric 69 and E the P6x0 text code - now D holds
#    the text code and E holds the numeric 69.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 11; $r++) {
    $eText = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 4).Value = $eText
    Set-NumericValue $ws.Cells.Item($r, 5) 69
}

# ---------------------------------------------------------------------
# 2. Insert a brand new row at position 12 for the "EindArtikel fallback"
#    entry - this pushes the old rows 12-15 down to 13-16 and the old
#    row 17 down to row 18, carrying their styles along with them.
# ---------------------------------------------------------------------
$ws.Rows.Item(12).Insert()

# Give the new row's D:J cells the same style (yellow fill, xf 4) as the
# rest of the "special" rows below it, matching A12/C12 which already
# inherited the plain style from row 11 above.
$ws.Range("D13:J13").Copy()
$ws.Range("D12:J12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A12").Value = "EindArtikel"
$ws.Range("C12").Value = "EindArtikel fallback (wanneer niet met een cijfer begint)"
$ws.Range("D12").Value = "P100"
Set-NumericValue $ws.Cells.Item(12, 5) 69
Set-NumericValue $ws.Cells.Item(12, 6) 6690
Set-NumericValue $ws.Cells.Item(12, 7) 86690
Set-NumericValue $ws.Cells.Item(12, 8) 30669
Set-NumericValue $ws.Cells.Item(12, 9) 76690
Set-NumericValue $ws.Cells.Item(12, 10) 81069

# ---------------------------------------------------------------------
# 3. Swap D/E columns for the rows that used to be 12-15 (now 13-16) too.
# ---------------------------------------------------------------------
for ($r = 13; $r -le 16; $r++) {
    $eText = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 4).Value = $eText
    Set-NumericValue $ws.Cells.Item($r, 5) 69
}

# Row 14 (old row 13's "GrondstofArtikel" line) had its J value switched
# from a text label to the plain number 81020.
Set-NumericValue $ws.Cells.Item(14, 10) 81020

# ---------------------------------------------------------------------
# 4. Sheet-level cosmetics: dimension / selection / column widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 28.140625
$ws.Columns.Item(5).ColumnWidth = 29

$ws.Range("A12").Select()
